# Update project summary with final product details:
# The syllabus entries (Topic / Do Before Class / In-Class Exercise, cols D:F)
# for the two class sessions on 2019-11-12 (row 15/16) and 2019-11-14/19
# (row 17/18) are swapped - i.e. the "Collaborating using Github" and
# "Big Data" sessions now come first (rows 15-16), followed by the
# "Pandas: Reshaping" and "Groupby" sessions (rows 17-18). Row heights
# (which are sized to fit the wrapped text) move along with the content.
# The active-cell selection also shifts one column right, from D19 to E19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch range well away from the used data to stage the swap so a
# direct D15:F16 <-> D17:F18 exchange doesn't clobber itself.
$scratch = $ws.Range("Z1:AB2")

# 1) Stash rows 15-16 (D:F) in the scratch area.
$ws.Range("D15:F16").Copy($scratch) | Out-Null

# 2) Move rows 17-18 (D:F) into rows 15-16.
$ws.Range("D17:F18").Copy($ws.Range("D15")) | Out-Null

# 3) Move the stashed original rows 15-16 into rows 17-18.
$scratch.Copy($ws.Range("D17")) | Out-Null

# 4) Clear the scratch area.
$scratch.Clear() | Out-Null

# 5) Swap the row heights that travel with the wrapped text so each row
#    is sized for its new content.
$h15 = $ws.Rows(15).RowHeight
$h16 = $ws.Rows(16).RowHeight
$h17 = $ws.Rows(17).RowHeight
$h18 = $ws.Rows(18).RowHeight
$ws.Rows(15).RowHeight = $h17
$ws.Rows(16).RowHeight = $h18
$ws.Rows(17).RowHeight = $h15
$ws.Rows(18).RowHeight = $h16

# 6) Move the active selection from D19 to E19.
$ws.Range("E19").Select() | Out-Null
